# Update in scripts 07/20
# Populate the AL_Policy Number and APD_Policy Number sheets with new
# policy number values (12 new rows on AL, 8 new rows on APD).

$wb = $excel.ActiveWorkbook

$wsAL = $wb.Worksheets.Item("AL_Policy Number")
$wsAPD = $wb.Worksheets.Item("APD_Policy Number")

# New Business (TUN-) policy numbers -> AL_Policy Number sheet, column A,
# starting at row 2 (row 1 holds the header).
$alValues = @(
    "TUN-OR-0001112",
    "TUN-OR-0001113",
    "TUN-MT-0001115",
    "TUN-CA-0001117",
    "TUN-IL-0001118",
    "TUN-IL-0001122",
    "TUN-IL-0001123",
    "TUN-IL-0001128",
    "TUN-IL-0001129",
    "TUN-IL-0001130",
    "TUN-MT-0001172",
    "TUN-MT-0001174"
)

# Out (TUO-) policy numbers -> APD_Policy Number sheet, column A,
# starting at row 2 (row 1 holds the header).
$apdValues = @(
    "TUO-OR-0001620",
    "TUO-MT-0001622",
    "TUO-CA-0001623",
    "TUO-IL-0001624",
    "TUO-IL-0001625",
    "TUO-IL-0001629",
    "TUO-MT-0001639",
    "TUO-MT-0001640"
)

for ($i = 0; $i -lt $alValues.Length; $i++) {
    $row = $i + 2
    $cell = $wsAL.Cells.Item($row, 1)
    $cell.Value = $alValues[$i]
    $cell.Style = "Normal"
}

for ($i = 0; $i -lt $apdValues.Length; $i++) {
    $row = $i + 2
    $cell = $wsAPD.Cells.Item($row, 1)
    $cell.Value = $apdValues[$i]
    $cell.Style = "Normal"
}
